# Update LR-pairs TPM data for Saa1-Fpr2 (Resolving-Mac replaces MuSCs as a sending cluster,
# and all dependent specificity / weight metrics are refreshed with the new TPM-derived values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9501423428580035
$ws.Range("J2").Value = 0.9501423428580037
$ws.Range("O2").Value = 0.2013489143128838
$ws.Range("P2").Value = 0.2013489143128838
$ws.Range("S2").Value = 0.1913101291771588
$ws.Range("T2").Value = 0.1913101291771588
$ws.Range("I3").Value = 0.9501423428580035
$ws.Range("J3").Value = 0.9501423428580037
$ws.Range("M3").Value = 0.246708
$ws.Range("N3").Value = 0.740124
$ws.Range("O3").Value = 0.03949536580856015
$ws.Range("P3").Value = 0.03949536580856015
$ws.Range("Q3").Value = 0.010374729288
$ws.Range("R3").Value = 0.093372563592
$ws.Range("S3").Value = 0.03752621940137922
$ws.Range("T3").Value = 0.03752621940137923
$ws.Range("I4").Value = 0.9501423428580035
$ws.Range("J4").Value = 0.9501423428580037
$ws.Range("M4").Value = 4.74207
$ws.Range("N4").Value = 14.22621
$ws.Range("O4").Value = 0.7591557198785561
$ws.Range("P4").Value = 0.759155719878556
$ws.Range("Q4").Value = 0.19941668902
$ws.Range("R4").Value = 1.79475020118
$ws.Range("S4").Value = 0.7213059942794655
$ws.Range("T4").Value = 0.7213059942794656
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002206666666666667
$ws.Range("H5").Value = 0.00662
$ws.Range("I5").Value = 0.04985765714199641
$ws.Range("J5").Value = 0.04985765714199642
$ws.Range("O5").Value = 0.2013489143128838
$ws.Range("P5").Value = 0.2013489143128838
$ws.Range("Q5").Value = 0.002775384246666667
$ws.Range("R5").Value = 0.02497845822
$ws.Range("S5").Value = 0.01003878513572497
$ws.Range("T5").Value = 0.01003878513572497
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.002206666666666667
$ws.Range("H6").Value = 0.00662
$ws.Range("I6").Value = 0.04985765714199641
$ws.Range("J6").Value = 0.04985765714199642
$ws.Range("M6").Value = 0.246708
$ws.Range("N6").Value = 0.740124
$ws.Range("O6").Value = 0.03949536580856015
$ws.Range("P6").Value = 0.03949536580856015
$ws.Range("Q6").Value = 0.0005444023200000001
$ws.Range("R6").Value = 0.00489962088
$ws.Range("S6").Value = 0.00196914640718092
$ws.Range("T6").Value = 0.00196914640718092
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.002206666666666667
$ws.Range("H7").Value = 0.00662
$ws.Range("I7").Value = 0.04985765714199641
$ws.Range("J7").Value = 0.04985765714199642
$ws.Range("M7").Value = 4.74207
$ws.Range("N7").Value = 14.22621
$ws.Range("O7").Value = 0.7591557198785561
$ws.Range("P7").Value = 0.759155719878556
$ws.Range("Q7").Value = 0.0104641678
$ws.Range("R7").Value = 0.0941775102
$ws.Range("S7").Value = 0.03784972559909052
$ws.Range("T7").Value = 0.03784972559909052
